# Add a new Job Posting row (Job_Id = JD_024) as the 25th data row (row 25)
# on the single worksheet. Mirrors: Job_Title="Mid Fullstack Developer",
# Job_Description = the standard "Junior RPA Developer" blurb (same text as
# used by rows 3, 4, 11, 19-22), Total_Years_Min_Exp=1, Total_Years_Max_Exp=3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "JD_024"
$ws.Range("B25").Value = "Mid Fullstack Developer"
$ws.Range("C25").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 3

# Re-fit the row height after entering the multi-line Job_Description so the
# row keeps the sheet's default (no leftover explicit/custom row height).
$ws.Rows.Item(25).EntireRow.AutoFit()
